$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the data rows (2-7). The rows are permuted as follows (old row -> new row):
#   2 -> 4, 3 -> 5, 4 -> 2, 5 -> 6, 6 -> 7, 7 -> 3
# Capture the original values first so we can rewrite them in the new order.
$orig = @{}
for ($r = 2; $r -le 7; $r++) {
    $orig[$r] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2
    )
}

$mapping = @{ 2 = 4; 3 = 5; 4 = 2; 5 = 6; 6 = 7; 7 = 3 }

foreach ($oldRow in $mapping.Keys) {
    $newRow = $mapping[$oldRow]
    $vals = $orig[$oldRow]
    $ws.Cells.Item($newRow, 1).Value = $vals[0]
    $ws.Cells.Item($newRow, 2).Value = $vals[1]
    $ws.Cells.Item($newRow, 3).Value = $vals[2]
}

# The special formatting that was on A6 (left/center aligned, Calibri Light)
# travels with its row content, which now lives in row 7.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active cell / selection as recorded in the saved view.
$ws.Range("F3").Select()
